$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Guide Quests")

# --- Update E2: remove <p><br></p> paragraphs between steps ---
$ws.Range("E2").Value = "<p>- Select a monster from the list behind this modal (start with Sewer Rat and work your way down)</p><p>- Click the attack button.</p><p>- Click the first attack button.</p><p>- Repeat until level 2.</p>"

# --- Update E3: remove stray trailing spaces and <p><br></p> paragraphs ---
$ws.Range("E3").Value = "<p>- First lets investigate the gear we have been given.</p><p>We can do this by going to the character sheet tab and on the bottom right is the inventory management. You can click item names here to investigate, equip and do other types of actions with the item. For now if you find gear that raises your stats (ie, STR Modifier +x%) equip it.</p><p>You may also need to visit the shop, which you can do by clicking the top left Hamburger menu to open the menu and select Shop. From here select Buy under General Shop. here you can buy gear, buy multiple pieces of gear or even compare and equip gear (auto purchase/equip).</p><p>- Next go back to the drop down for monsters and select a stronger monster and click Attack</p><p>- Select Attack and if you can kill it in one hit, which is ideal, move down the list to the next, rinse and repeat till you cannot move any further.</p><p>- Now that we have a monster, click Explore to the left of the attack section.</p><p>- Select the same monster, select 1 hour, ignore the move down and then select Attack, click Explore.</p><p>This will run a set of battles every 5 minutes for 1 hour where you can fight between 1-8 enemies back to back. The reason we ignored the move down aspect is because it allows you to state: Move down the list every x levels that I gain, and for now we do not want that.</p>"

# Column E's width is a "best fit" width that Excel recalculates whenever the
# underlying text changes length; reproduce the narrower fit for the shorter text.
$ws.Columns.Item(5).ColumnWidth = 1594

# --- Decrement the id values in column A for rows 5 through 15 ---
for ($r = 5; $r -le 15; $r++) {
    $current = [double]$ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $current - 1
}
